$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text-typed cell value without Excel auto-converting
# numeric- or percent-looking strings into Number/Percentage cells.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# row, column letter, new value
$updates = @(
    @(2, "D", "304.81"),
    @(2, "E", "4.00%"),
    @(3, "D", "32.34"),
    @(3, "E", "5.63%"),
    @(4, "E", "2.85%"),
    @(5, "D", "0.07615"),
    @(5, "E", "6.82%"),
    @(6, "B", "FTXToken"),
    @(6, "C", "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"),
    @(6, "D", "1.874"),
    @(6, "E", "33.32%"),
    @(7, "B", "KuCoinToken"),
    @(7, "C", "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"),
    @(7, "D", "7.894"),
    @(7, "E", "4.53%"),
    @(8, "B", "GateToken"),
    @(8, "C", "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"),
    @(8, "D", "3.876"),
    @(8, "E", "7.21%"),
    @(9, "D", "0.9310"),
    @(9, "E", "1.56%"),
    @(10, "D", "0.1713"),
    @(10, "E", "4.35%"),
    @(11, "D", "0.07980"),
    @(11, "E", "0.65%"),
    @(12, "D", "0.08043"),
    @(12, "E", "3.34%"),
    @(13, "D", "0.03052"),
    @(13, "E", "3.54%"),
    @(14, "D", "0.09934"),
    @(14, "E", "10.37%"),
    @(15, "D", "0.001490"),
    @(15, "E", "-5.50%"),
    @(16, "D", "0.04601"),
    @(16, "E", "1.35%"),
    @(17, "D", "0.006304"),
    @(17, "E", "1.15%"),
    @(19, "D", "2.231"),
    @(19, "E", "-0.57%"),
    @(20, "E", "1.50%"),
    @(21, "D", "0.1343"),
    @(21, "E", "-0.18%"),
    @(22, "D", "4.546"),
    @(22, "E", "9.60%"),
    @(23, "E", "1.64%"),
    @(24, "D", "0.001215"),
    @(24, "E", "0.49%"),
    @(25, "D", "0.004501"),
    @(25, "E", "6.16%"),
    @(26, "E", "19.50%"),
    @(27, "D", "0.0001780"),
    @(27, "E", "5.43%"),
    @(39, "D", "0.01738"),
    @(39, "E", "2,545.72%"),
    @(40, "D", "0.04547"),
    @(40, "E", "2.84%"),
    @(41, "D", "0.006967"),
    @(41, "E", "-1.16%"),
    @(42, "D", "0.1362"),
    @(42, "E", "6.90%"),
    @(43, "D", "0.01390"),
    @(43, "E", "5.11%"),
    @(44, "D", "0.002066"),
    @(44, "E", "-5.60%"),
    @(45, "D", "0.00006139"),
    @(45, "E", "4.92%"),
    @(46, "D", "0.7091"),
    @(46, "E", "-62.73%"),
    @(47, "D", "0.01220"),
    @(47, "E", "-6.04%")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $newValue = $u[2]
    $colIndex = switch ($col) {
        "B" { 2 }
        "C" { 3 }
        "D" { 4 }
        "E" { 5 }
    }
    if ($col -eq "D" -or $col -eq "E") {
        Set-TextValue $row $colIndex $newValue
    } else {
        $ws.Cells.Item($row, $colIndex).Value = $newValue
    }
}
